# Update odds data for a handful of matches (rows 5, 6, 8, 15) in the
# "Jogos da Semana" FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 - ANTIGUA & BARBUDA - ABFA PREMIER LEAGUE | Five Islands - Old Road
$ws.Range("G5").Value = 24
$ws.Range("H5").Value = 9.25
$ws.Range("I5").Value = 1.04
$ws.Range("J5").Value = 14.5
$ws.Range("K5").Value = 4.05
$ws.Range("L5").Value = 1.23
$ws.Range("S5").Value = 1.04
$ws.Range("T5").Value = 8.4
$ws.Range("U5").Value = 1.69
$ws.Range("V5").Value = 2.11
$ws.Range("W5").Value = 200
$ws.Range("Y5").Value = 120
$ws.Range("AA5").Value = 500
$ws.Range("AB5").Value = 175
$ws.Range("AC5").Value = 60
$ws.Range("AD5").Value = 35
$ws.Range("AE5").Value = 37
$ws.Range("AF5").Value = 80
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 16
$ws.Range("AK5").Value = 9
$ws.Range("AL5").Value = 11
$ws.Range("AM5").Value = 25
$ws.Range("AN5").Value = 29
$ws.Range("AO5").Value = 175
$ws.Range("AP5").Value = 65
$ws.Range("AR5").Value = 500
$ws.Range("AS5").Value = 400
$ws.Range("AT5").Value = 7.8
$ws.Range("AU5").Value = 11
$ws.Range("AV5").Value = 50
$ws.Range("AW5").Value = 500
$ws.Range("AX5").Value = 4.25
$ws.Range("AY5").Value = 4.3
$ws.Range("AZ5").Value = 11
$ws.Range("BA5").Value = 6.3
$ws.Range("BB5").Value = 16
$ws.Range("BC5").Value = 75

# Row 6 - ARGENTINA - TORNEO BETANO | Central Cordoba - Rosario Central
$ws.Range("Q6").Value = 3.1
$ws.Range("R6").Value = 1.36

# Row 8 - ARGENTINA - TORNEO BETANO | Belgrano - Ind. Rivadavia
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 2.63
$ws.Range("L8").Value = 5
$ws.Range("S8").Value = 1.53
$ws.Range("T8").Value = 2.38
$ws.Range("U8").Value = 2.2
$ws.Range("V8").Value = 1.62
$ws.Range("X8").Value = 8
$ws.Range("Y8").Value = 9
$ws.Range("AC8").Value = 7
$ws.Range("AJ8").Value = 15
$ws.Range("AN8").Value = 3.75
$ws.Range("AP8").Value = 26
$ws.Range("AS8").Value = 251
$ws.Range("AT8").Value = 2.38
$ws.Range("AY8").Value = 26
$ws.Range("BC8").Value = 401
$ws.Range("BD8").Value = 151

# Row 15 - BOLIVIA - DIVISION PROFESIONAL | Oriente Petrolero - Santa Cruz
$ws.Range("G15").Value = 1.45
$ws.Range("I15").Value = 7.5
$ws.Range("K15").Value = 2.38
$ws.Range("AE15").Value = 19
$ws.Range("AG15").Value = 351
$ws.Range("AO15").Value = 7
